$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("06-10-2021", "07-10-2021", "08-10-2021", "12-10-2021")
$row = 193
foreach ($d in $dates) {
    $cellA = $ws.Cells.Item($row, 1)
    # Enter as a text-returning formula first so Excel does not
    # auto-convert the dd-mm-yyyy-looking string into a date serial,
    # then convert the formula result to a static value (paste values)
    # so the final cell is a plain shared-string cell like its neighbours.
    $cellA.Formula = "=""" + $d + """"
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    $ws.Cells.Item($row, 2).Value = 3.25

    $row = $row + 1
}
$excel.CutCopyMode = $false
